$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Insert a new column at S, shifting the existing S:AB headers one column to the right
$ws.Range("S1").EntireColumn.Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)

# New header cell: text + formatting to match the other header cells in row 1
# (bold font, left/center alignment with indent - same direct formatting as style index 2)
$ws.Range("S1").Value = "Account type"
$ws.Range("S1").Font.Bold = $true
$ws.Range("S1").HorizontalAlignment = -4131
$ws.Range("S1").VerticalAlignment = -4108
$ws.Range("S1").IndentLevel = 1
$ws.Range("S1").ColumnWidth = 18.2

# Restore the previous selection / scroll position used on this sheet
$ws.Activate()
$ws.Range("Q7").Select()
